$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Certificados, cursos, badges")

# New row of data (row 73), mirroring the formatting of the previous row (72)
$row = 73
$prevRow = 72

$ws.Range("B$prevRow`:I$prevRow").Copy() | Out-Null
$ws.Range("B$row`:I$row").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item($row, 2).Value = "Data Science Academy"
$ws.Cells.Item($row, 3).Value = "Engenharia de Dados com Hadoop e Spark"
$ws.Cells.Item($row, 4).Value = 64
$ws.Cells.Item($row, 5).Value = "08/30/2024"

$ws.Hyperlinks.Add($ws.Range("F$row"), "https://mycourse.app/BLDUSreYgR2wGoaKA") | Out-Null

# Hyperlinks.Add resets the cell's font to the default Hyperlink style;
# reapply the same formatting used across the rest of column F.
$ws.Range("F$prevRow").Copy() | Out-Null
$ws.Range("F$row").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item($row, 7).Value = "Ok"
$ws.Cells.Item($row, 8).Value = "Ok"
$ws.Cells.Item($row, 9).Value = "08/30/2024"

# Update the view so row 73 is visible / selected, similar to final workbook state
$ws.Range("B73").Select()
$excel.ActiveWindow.ScrollRow = 36
